$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Insert 3 new rows at row 13 or the two extra "Docentes responsaveis"
# names (pushes the old rows 13-24 down to rows 16-27).
# ---------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# ---------------------------------------------------------------
# Row 10/11 (Objetivos:/Objectives:) - fill in the real text (previously
# these cells erroneously held a professor name)
# ---------------------------------------------------------------
$ws.Range('B10').Value = 'A disciplina Degradação e Proteção de Materiais visa propiciar aos alunos os conhecimentos básicos da degradação dos materiais, dando enfoque à corrosão e à oxidação de metais e ligas metálicas. São abordados os fundamentos teóricos dos dois processos degradativos, tanto termodinâmicos como cinéticos, e descritas as principais formas de ataque e as técnicas de proteção contra a corrosão e a oxidação metálica.'
$ws.Range('C10').Value = 'A disciplina Degradação e Proteção de Materiais visa propiciar aos alunos os conhecimentos básicos da degradação dos materiais, dando enfoque à corrosão e à oxidação de metais e ligas metálicas. São abordados os fundamentos teóricos dos dois processos degradativos, tanto termodinâmicos como cinéticos, e descritas as principais formas de ataque e as técnicas de proteção contra a corrosão e a oxidação metálica.'

# ---------------------------------------------------------------
# New rows 13-15: the three "Docentes responsaveis" entries.
# Copy B3:C3 formatting (style 2 / style 3, no custom row height) onto
# them first, then set the values.
# ---------------------------------------------------------------
$ws.Range('B3:C3').Copy()
$ws.Range('B13:C15').PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range('B13').Value = '5840747 - Alain Laurent Marie Robin'
$ws.Range('C13').Value = '5840747 - Alain Laurent Marie Robin'
$ws.Range('B14').Value = '7926291 - Célia Regina Tomachuk dos Santos Catuogno'
$ws.Range('C14').Value = '7926291 - Célia Regina Tomachuk dos Santos Catuogno'
$ws.Range('B15').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range('C15').Value = '7459752 - Maria Ismenia Sodero Toledo Faria'

# ---------------------------------------------------------------
# Remaining content updates on the rows that shifted down (old rows
# 13,14,15,18,19,20,21 -> new rows 16,17,18,21,22,23,24). Only the
# B/C text content changes; formatting & row heights already carried
# over with the row shift.
# ---------------------------------------------------------------
$ws.Range('B16').Value = 'Fundamentos da corrosão (termodinâmica e cinética); principais tipos de corrosão; Controle e proteção contra a corrosão; Degradação de materiais poliméricos e cerâmicos; Oxidação em altas temperaturas. Estudo de Caso'
$ws.Range('C16').Value = 'Fundamentos da corrosão (termodinâmica e cinética); principais tipos de corrosão; Controle e proteção contra a corrosão; Degradação de materiais poliméricos e cerâmicos; Oxidação em altas temperaturas. Estudo de Caso'
$ws.Range('B18').Value = 'Corrosão e sua importância econômica, social e ambiental. 2. Aspectos termodinâmicos e cinéticos da corrosão. Polarização. Princípios básicos de eletroquímica para compreensão do fenômeno de corrosão. Diagrama de Pourbaix. 3. Tipos de corrosão (corrosão uniforme, corrosão por pites, corrosão intergranular, corrosão associado com fatores mecânicos, corrosão galvânica e corrosão atmosférica). 4. Proteção catódica e anódica. 5.Revestimentos metálicos e orgânicos. 6. Inibidores de corrosão. 7. Degradação de materiais poliméricos. 8.Degradação de materiais cerâmicos. 9.Oxidação em altas temperaturas. 10. Estudos de Casos referentes às falhas causadas em equipamentos devido processo corrosivo, atividade que potencializa o papel do aluno como protagonista do processo de ensino e aprendizagem, colocando-o em contato com problemas reais.'
$ws.Range('C18').Value = 'Corrosão e sua importância econômica, social e ambiental. 2. Aspectos termodinâmicos e cinéticos da corrosão. Polarização. Princípios básicos de eletroquímica para compreensão do fenômeno de corrosão. Diagrama de Pourbaix. 3. Tipos de corrosão (corrosão uniforme, corrosão por pites, corrosão intergranular, corrosão associado com fatores mecânicos, corrosão galvânica e corrosão atmosférica). 4. Proteção catódica e anódica. 5.Revestimentos metálicos e orgânicos. 6. Inibidores de corrosão. 7. Degradação de materiais poliméricos. 8.Degradação de materiais cerâmicos. 9.Oxidação em altas temperaturas. 10. Estudos de Casos referentes às falhas causadas em equipamentos devido processo corrosivo, atividade que potencializa o papel do aluno como protagonista do processo de ensino e aprendizagem, colocando-o em contato com problemas reais.'
$ws.Range('B21').Value = 'Os alunos serão avaliados continuamente quanto às habilidades gerais em função da participação ativa nas aulas. Também, serão aplicadas: provas escritas, trabalhos extraclasse, pequenos seminários e atividades para discussão dos Estudos de Caso.'
$ws.Range('C21').Value = 'Os alunos serão avaliados continuamente quanto às habilidades gerais em função da participação ativa nas aulas. Também, serão aplicadas: provas escritas, trabalhos extraclasse, pequenos seminários e atividades para discussão dos Estudos de Caso.'
$ws.Range('B22').Value = 'Será considerada a média das avaliações gerais, com peso 1 (AG) somada à nota de duas provass escritas, P1 e P2, com peso 1, cada). A nota final (NF) será calculada pela equação: NF = (AG + P1+ P2)/3.'
$ws.Range('C22').Value = 'Será considerada a média das avaliações gerais, com peso 1 (AG) somada à nota de duas provass escritas, P1 e P2, com peso 1, cada). A nota final (NF) será calculada pela equação: NF = (AG + P1+ P2)/3.'
$ws.Range('B23').Value = 'Prova escrita sobre toda matéria.A média final MF será a média da nota final NF e da nota obtida na recuperação NR:MF = (NF + NR)/2Será aprovado o aluno com MF igual ou superior a 5.'
$ws.Range('C23').Value = 'Prova escrita sobre toda matéria.A média final MF será a média da nota final NF e da nota obtida na recuperação NR:MF = (NF + NR)/2Será aprovado o aluno com MF igual ou superior a 5.'
$ws.Range('B24').Value = 'ASM Handbook Committee, ASM Handbook: Corrosion: Fundamentals, Testing, and Protection, vol. 13A, Materials Park, Ohio: ASM International, 2003.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Coimbra, Livraria Medina, 1996. FONTANA, M.G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987  GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro,Ed. LTC, 2007 JAMBO, H.C.M. e Fofano S. Corrosão: Fundamentos, Monitoração e Controle. Editora Ciência Moderna,2009. JONES, D.A. Principles and Prevention of Corrosion. 2ª Edição, Prentice Hall, 1996. McCAULEY, R. A. Corrosion of Ceramic and Composite Materials. 2ª Edição, 2004. RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus,1990. SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2,2000. WEST, J. M. Basic Corrosion and Oxidation. Second Edition. New York. Ellis Horwood Limited; John Wiley & Sons.1986.ROBERGE,Pierre R. Corrosion Engineering Principles and PracticeMcGraw-Hill, 2008..'
$ws.Range('C24').Value = 'ASM Handbook Committee, ASM Handbook: Corrosion: Fundamentals, Testing, and Protection, vol. 13A, Materials Park, Ohio: ASM International, 2003.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Coimbra, Livraria Medina, 1996. FONTANA, M.G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987  GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro,Ed. LTC, 2007 JAMBO, H.C.M. e Fofano S. Corrosão: Fundamentos, Monitoração e Controle. Editora Ciência Moderna,2009. JONES, D.A. Principles and Prevention of Corrosion. 2ª Edição, Prentice Hall, 1996. McCAULEY, R. A. Corrosion of Ceramic and Composite Materials. 2ª Edição, 2004. RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus,1990. SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2,2000. WEST, J. M. Basic Corrosion and Oxidation. Second Edition. New York. Ellis Horwood Limited; John Wiley & Sons.1986.ROBERGE,Pierre R. Corrosion Engineering Principles and PracticeMcGraw-Hill, 2008..'

Write-Output "done"